# Rescale the clustering counts and incorporate the "time" dimension:
# the original 2-row table (labels 0,1) becomes a 4-row table (labels 3,0,1,2)
# with updated counts, so the sheet grows from A1:B3 to A1:B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update / add the data rows (row 1 with B1=0 stays untouched)
$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(2, 2).Value = 94

$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 2).Value = 92

$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = 71

$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = 58

# The new rows 4 and 5 in column A need the same formatting (style) that
# column A already carries in rows 2-3, so copy it down.
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
